# Updated cryptos list with GitHub Actions: refresh Price/Volume(1h) columns,
# and re-rank two coin pairs that swapped positions (Bittensor/PEPE and USDe/InjectiveProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, since several "Price" values
# look numeric (e.g. "608.34", "35.20", "2.30") and must stay literal strings
# (matching the sheet's existing inlineStr cells) rather than being coerced
# into floating point numbers that would lose trailing zeros / precision.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.960.33"
$ws.Range("E2").Value = "  -4.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.138.14"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.34"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.45"
$ws.Range("E6").Value = "  -7.22%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.130.45"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("E9").Value = "  -5.10%  "
$ws.Range("E10").Value = "  -7.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("E11").Value = "  -7.80%  "
$ws.Range("E12").Value = "  -6.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -8.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.20"
$ws.Range("E14").Value = "  -10.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.661.02"
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.998.19"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.142.88"
$ws.Range("E18").Value = "  -3.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  -8.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.52"
$ws.Range("E20").Value = "  -6.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.68"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -6.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.74"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.53"
$ws.Range("E24").Value = "  -8.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.61"
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.37"
$ws.Range("E28").Value = "  -8.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -10.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.71"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("E31").Value = "  -15.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.99"
$ws.Range("E34").Value = "  -7.54%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  -8.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.52"
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "460.70"
$ws.Range("E38").Value = "  -7.11%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0728"
$ws.Range("E39").Value = "  -7.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  -14.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0392"
$ws.Range("E41").Value = "  -8.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.37"
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  -8.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.835.04"
$ws.Range("E44").Value = "  -5.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("E45").Value = "  -10.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -10.95%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.26"
$ws.Range("E48").Value = "  -9.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  -8.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.96"
$ws.Range("E51").Value = "  -1.78%  "
